$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") corrections
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total") corrections
$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -14
$ws.Range("E12").Value = "22 / 112"
